# petty-cashBook-2021.xlsx -- "Update 3-Jun-2021, midday update."
# Adds 7 new ledger entries (rows 19-25) on the 44349 (2-Jun-2021) date
# block of "Buku KAS HARIAN"-style Sheet1, a new date row (26, 3-Jun-2021),
# and moves the viewport/selection further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 19: Wages Expense ------------------------------------------------
$ws.Range("B19").Value = "Wages Expense"
$ws.Range("D19").Formula = "=60000+300000"

# --- Row 20: TRANSFER BCA --------------------------------------------------
$ws.Range("B20").Value = "TRANSFER BCA"
$ws.Range("D20").Formula = "=611500+385000+4557500+119000+5410000+16180000+521000+40000+4025000"

# --- Row 21: A/R ------------------------------------------------------------
$ws.Range("B21").Value = "A/R"
$ws.Range("C21").Formula = "=4557500+119000+2700000+31598500"

# --- Row 22: GARRETH - buku piano (new payee) -------------------------------
$ws.Range("B22").Value = "GARRETH - buku piano"
$ws.Range("D22").Value = 515000

# --- Row 23: BENSIN - RUSH (new payee) --------------------------------------
$ws.Range("B23").Value = "BENSIN - RUSH"
$ws.Range("D23").Value = 250000

# --- Row 24: SALES - cash/retail --------------------------------------------
$ws.Range("B24").Value = "SALES - cash/retail"
$ws.Range("C24").Formula = "=24619375+18942625-31598500"

# --- Row 25: SETOR KE BANK --------------------------------------------------
$ws.Range("B25").Value = "SETOR KE BANK"
$ws.Range("D25").Value = 18000000

# --- Row 26: new date (3-Jun-2021 = serial 44350) ---------------------------
$ws.Range("A26").Value = 44350

# --- Viewport: scroll the frozen pane further down, move the selection -----
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("C46").Select()
